$d = $word.ActiveDocument

# The block being removed is the whole "pricing pseudocode" section that
# starts with the paragraph "IF beauty = (plaiting) THEN " and runs
# through the paragraph "READ total" (inclusive). Locate both ends by
# their paragraph text so the edit is resilient to any paragraph-index
# drift, then delete the Range spanning them.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t -like "IF beauty*THEN*") {
        $startPara = $p
    }
    if ($t -like "READ total*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
